# This script updates the "dSF" column (F) values for a set of rows in
# Sheet1, reflecting a repull/recalculation of the underlying data.
# (commit message: "repull data, push all data, mean calculation")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = -2
    8  = 3
    11 = 3
    15 = -2
    16 = 2
    17 = -3
    20 = -1
    23 = -5
    24 = -2
    28 = -7
    35 = -4
    36 = -9
    42 = -1
    47 = -3
    49 = -4
    54 = -12
    59 = -2
    60 = -3
    65 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
